# Reorder the "Recorded By" (column G) entries so that the literal token
# "System" (exact case) is moved to the front of the comma-separated list,
# while preserving the original relative order of all remaining tokens.
# Rows whose G value does not contain an exact "System" token, or that
# consist of "System" alone, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G$row")
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $parts = $val.Split(",")

    if ($parts.Length -le 1) {
        continue
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Trim().Equals("System")) {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        continue
    }

    $rest = @()
    foreach ($p in $parts) {
        $t = $p.Trim()
        if (-not $t.Equals("System")) {
            $rest += $t
        }
    }

    $newVal = "System, " + ($rest -join ", ")

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
